$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo "tody" -> "today" in the vaccination question (appears in column J, row 56) ---
$ws.Range("J56").Value = "How many people were vaccinated today?"

# --- Column A (curated/easy questions) is updated:
#       * typo "tody" -> "today" fixed
#       * "how many cases do we have in austria" removed
#       * "Which countries reported on the 7th of February 2022 more than 100.000 new cases?" removed
#       * "How you doin'?" appended at the end
#     Rows 1-37 of column A are untouched; rows 38-45 get the rewritten tail of the list,
#     and the now-unused rows 46-49 are cleared out entirely.
$ws.Range("A38").Value = "How many people are vaccinated?"
$ws.Range("A39").Value = "How many people were vaccinated today?"
$ws.Range("A40").Value = "How many corona cases got reported today?"
$ws.Range("A41").Value = "when did austria have the highest number of infections?"
$ws.Range("A42").Value = "On which day was the highest number of new cases reported in Brazil?"
$ws.Range("A43").Value = "How many vaccine shots were administered in Austria on the 9th February 2022?"
$ws.Range("A44").Value = "Could you tell me which how many new COVID cases were reported in Austria on 2nd February 2022?"
$ws.Range("A45").Value = "How you doin'?"
$ws.Range("A46:A49").Clear()

# The whole populated column A range (A1:A45) uses the "Good" (green) cell style.
$ws.Range("A1:A45").Style = "Good"

# --- Highlight J32 with the new "Neutral" (yellow) cell style and bump the row height ---
$ws.Range("J32").Style = "Neutral"
$ws.Rows(32).RowHeight = 16

# --- Update the saved view/selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A45").Select()
